$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '294.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.49%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.14'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.00%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.97%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07341'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.09%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.287'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '30.06%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.743'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.89%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.30%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9085'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.54%'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.96%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08035'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '7.28%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08062'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.28%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03101'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.84%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.86%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001516'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.71%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005719'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.52%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.477'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.50%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.075'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.33%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3328'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.53%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.38%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.980'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-9.76%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.85%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04548'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.81%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001211'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.17%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '15.58%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.93%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003398'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01602'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-2.41%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04451'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.74%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007319'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.34%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1331'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008628'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.001948'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.43%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.25%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00005971'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '4.30%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.06%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '2.38%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-3.52%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.06%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.06%'
